$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows are appended as plain text (matching the existing sheet's
# convention where every cell - including numeric-looking ones - is stored
# as text), so force a Text number format before writing numeric-looking
# values to stop Excel from auto-coercing them into real numbers. The style
# is reset back to "Normal" afterwards so the cells keep the workbook's
# original (unstyled) look, same as every other cell on the sheet.
$ws.Range("A4:K5").NumberFormat = "@"

# Row 4: " Abu Dhabi" / Mumbai match (same details as the existing row 3)
$ws.Cells.Item(4, 1).Value = " Abu Dhabi"
$ws.Cells.Item(4, 2).Value = " September 23 2020"
$ws.Cells.Item(4, 3).Value = "Mumbai won by 49 runs"
$ws.Cells.Item(4, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(4, 5).Value = "Mumbai Indians"
$ws.Cells.Item(4, 6).Value = "Kuldeep Yadav "
$ws.Cells.Item(4, 7).Value = "1"
$ws.Cells.Item(4, 8).Value = "2"
$ws.Cells.Item(4, 9).Value = "0"
$ws.Cells.Item(4, 10).Value = "0"
$ws.Cells.Item(4, 11).Value = "50.00"

# Row 5: " Abu Dhabi" / RCB match (same details as the existing row 2)
$ws.Cells.Item(5, 1).Value = " Abu Dhabi"
$ws.Cells.Item(5, 2).Value = " October 21 2020"
$ws.Cells.Item(5, 3).Value = "RCB won by 8 wickets (with 39 balls remaining)"
$ws.Cells.Item(5, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(5, 5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(5, 6).Value = "Kuldeep Yadav "
$ws.Cells.Item(5, 7).Value = "12"
$ws.Cells.Item(5, 8).Value = "19"
$ws.Cells.Item(5, 9).Value = "1"
$ws.Cells.Item(5, 10).Value = "0"
$ws.Cells.Item(5, 11).Value = "63.15"

# Drop the temporary Text formatting so the new cells end up with the same
# (default) style as the rest of the sheet.
$ws.Range("A4:K5").Style = "Normal"
